$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
  @{row=2; A="ECs"; B="Lrrc4b"; C="Ptprf"; D="ECs"; E=1; F=0.3333333333333333; G=0.063665; H=0.190995; I=0.2699335606864689; J=0.269933560686469; K=3; L=1; M=0.4211356666666666; N=1.263407; O=0.05962259118326733; P=0.05962259118326733; Q=0.02681160221833333; R=0.241304419965; S=0.01609413833545302; T=0.01609413833545302}
  @{row=3; A="ECs"; B="Lrrc4b"; C="Ptprf"; D="FAPs"; E=1; F=0.3333333333333333; G=0.063665; H=0.190995; I=0.2699335606864689; J=0.269933560686469; K=3; L=1; M=5.467658; N=16.402974; O=0.7740876954075475; P=0.7740876954075474; Q=0.34809844657; R=3.13288601913; S=0.2089522479049421; T=0.2089522479049421}
  @{row=4; A="ECs"; B="Lrrc4b"; C="Ptprf"; D="Inflammatory-Mac"; E=1; F=0.3333333333333333; G=0.063665; H=0.190995; I=0.2699335606864689; J=0.269933560686469; K=3; L=1; M=0.3825986666666667; N=1.147796; O=0.05416668711649494; P=0.05416668711649493; Q=0.02435814411333333; R=0.21922329702; S=0.01462140672394536; T=0.01462140672394536}
  @{row=5; A="ECs"; B="Lrrc4b"; C="Ptprf"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.063665; H=0.190995; I=0.2699335606864689; J=0.269933560686469; K=3; L=1; M=0.7919649999999999; N=2.375895; O=0.1121230262926903; P=0.1121230262926903; Q=0.050420451725; R=0.4537840655249999; S=0.03026576772212845; T=0.03026576772212846}
  @{row=6; A="FAPs"; B="Lrrc4b"; C="Ptprf"; D="ECs"; E=1; F=0.3333333333333333; G=0.1436256666666667; H=0.430877; I=0.6089592022194491; J=0.6089592022194491; K=3; L=1; M=0.4211356666666666; N=1.263407; O=0.05962259118326733; P=0.05962259118326733; Q=0.06048589088211111; R=0.544373017939; S=0.03630772556121883; T=0.03630772556121883}
  @{row=7; A="FAPs"; B="Lrrc4b"; C="Ptprf"; D="FAPs"; E=1; F=0.3333333333333333; G=0.1436256666666667; H=0.430877; I=0.6089592022194491; J=0.6089592022194491; K=3; L=1; M=5.467658; N=16.402974; O=0.7740876954075475; P=0.7740876954075474; Q=0.7852960253553334; R=7.067664228198001; S=0.4713878254432721; T=0.471387825443272}
  @{row=8; A="FAPs"; B="Lrrc4b"; C="Ptprf"; D="Inflammatory-Mac"; E=1; F=0.3333333333333333; G=0.1436256666666667; H=0.430877; I=0.6089592022194491; J=0.6089592022194491; K=3; L=1; M=0.3825986666666667; N=1.147796; O=0.05416668711649494; P=0.05416668711649493; Q=0.05495098856577778; R=0.494558897092; S=0.03298530257333127; T=0.03298530257333127}
  @{row=9; A="FAPs"; B="Lrrc4b"; C="Ptprf"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.1436256666666667; H=0.430877; I=0.6089592022194491; J=0.6089592022194491; K=3; L=1; M=0.7919649999999999; N=2.375895; O=0.1121230262926903; P=0.1121230262926903; Q=0.1137465011016667; R=1.023718509915; S=0.06827834864162698; T=0.06827834864162698}
  @{row=10; A="Inflammatory-Mac"; B="Lrrc4b"; C="Ptprf"; D="ECs"; E=1; F=0.3333333333333333; G=0.02856366666666667; H=0.085691; I=0.1211072370940821; J=0.1211072370940821; K=3; L=1; M=0.4211356666666666; N=1.263407; O=0.05962259118326733; P=0.05962259118326733; Q=0.01202917880411111; R=0.108262609237; S=0.007220727286595484; T=0.007220727286595485}
  @{row=11; A="Inflammatory-Mac"; B="Lrrc4b"; C="Ptprf"; D="FAPs"; E=1; F=0.3333333333333333; G=0.02856366666666667; H=0.085691; I=0.1211072370940821; J=0.1211072370940821; K=3; L=1; M=5.467658; N=16.402974; O=0.7740876954075475; P=0.7740876954075474; Q=0.1561763605593333; R=1.405587245034; S=0.09374762205933346; T=0.09374762205933344}
  @{row=12; A="Inflammatory-Mac"; B="Lrrc4b"; C="Ptprf"; D="Inflammatory-Mac"; E=1; F=0.3333333333333333; G=0.02856366666666667; H=0.085691; I=0.1211072370940821; J=0.1211072370940821; K=3; L=1; M=0.3825986666666667; N=1.147796; O=0.05416668711649494; P=0.05416668711649493; Q=0.01092842078177778; R=0.09835578703600001; S=0.006559977819218314; T=0.006559977819218314}
  @{row=13; A="Inflammatory-Mac"; B="Lrrc4b"; C="Ptprf"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.02856366666666667; H=0.085691; I=0.1211072370940821; J=0.1211072370940821; K=3; L=1; M=0.7919649999999999; N=2.375895; O=0.1121230262926903; P=0.1121230262926903; Q=0.02262142427166667; R=0.203592818445; S=0.01357890992893484; T=0.01357890992893484}
)

foreach ($item in $rowsData) {
  $r = $item.row
  if ($item.ContainsKey("A")) { $ws.Range("A$r").Value = $item.A }
  if ($item.ContainsKey("B")) { $ws.Range("B$r").Value = $item.B }
  if ($item.ContainsKey("C")) { $ws.Range("C$r").Value = $item.C }
  if ($item.ContainsKey("D")) { $ws.Range("D$r").Value = $item.D }
  if ($item.ContainsKey("E")) { $ws.Range("E$r").Value = $item.E }
  if ($item.ContainsKey("F")) { $ws.Range("F$r").Value = $item.F }
  if ($item.ContainsKey("G")) { $ws.Range("G$r").Value = $item.G }
  if ($item.ContainsKey("H")) { $ws.Range("H$r").Value = $item.H }
  if ($item.ContainsKey("I")) { $ws.Range("I$r").Value = $item.I }
  if ($item.ContainsKey("J")) { $ws.Range("J$r").Value = $item.J }
  if ($item.ContainsKey("K")) { $ws.Range("K$r").Value = $item.K }
  if ($item.ContainsKey("L")) { $ws.Range("L$r").Value = $item.L }
  if ($item.ContainsKey("M")) { $ws.Range("M$r").Value = $item.M }
  if ($item.ContainsKey("N")) { $ws.Range("N$r").Value = $item.N }
  if ($item.ContainsKey("O")) { $ws.Range("O$r").Value = $item.O }
  if ($item.ContainsKey("P")) { $ws.Range("P$r").Value = $item.P }
  if ($item.ContainsKey("Q")) { $ws.Range("Q$r").Value = $item.Q }
  if ($item.ContainsKey("R")) { $ws.Range("R$r").Value = $item.R }
  if ($item.ContainsKey("S")) { $ws.Range("S$r").Value = $item.S }
  if ($item.ContainsKey("T")) { $ws.Range("T$r").Value = $item.T }
}
